# Updated cryptos list on Fri Jun 28 04:54:05 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# snapshot table on the active sheet with the latest scraped figures.
# Numeric-looking Price values are entered with a leading apostrophe so
# Excel keeps them as text (matching the sheet's existing text-formatted
# price strings, e.g. "61.627.92") instead of silently coercing them into
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.627.92'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '3.452.47'
$ws.Range("E3").Value = '  +2.06%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'581.26"
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("D6").Value = "'145.97"
$ws.Range("E6").Value = '  +6.49%  '
$ws.Range("D7").Value = '3.454.18'
$ws.Range("E7").Value = '  +2.13%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = "'0.477"
$ws.Range("E9").Value = '  +1.58%  '
$ws.Range("E10").Value = '  +0.12%  '
$ws.Range("E11").Value = '  +2.76%  '
$ws.Range("E12").Value = '  +2.47%  '
$ws.Range("D13").Value = '4.041.58'
$ws.Range("E13").Value = '  +2.10%  '
$ws.Range("D14").Value = "'28.02"
$ws.Range("E14").Value = '  +8.88%  '
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").Value = '3.455.42'
$ws.Range("E17").Value = '  +2.11%  '
$ws.Range("D18").Value = '61.751.45'
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("E19").Value = '  +8.51%  '
$ws.Range("D20").Value = "'14.36"
$ws.Range("E20").Value = '  +3.91%  '
$ws.Range("D21").Value = "'9.55"
$ws.Range("E21").Value = '  +2.17%  '
$ws.Range("D22").Value = "'390.92"
$ws.Range("E22").Value = '  +3.89%  '
$ws.Range("E23").Value = '  +2.98%  '
$ws.Range("D24").Value = "'73.81"
$ws.Range("E24").Value = '  +3.89%  '
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("E27").Value = '  -0.97%  '
$ws.Range("D28").Value = '3.592.76'
$ws.Range("E28").Value = '  +2.08%  '
$ws.Range("D29").Value = "'0.183"
$ws.Range("E29").Value = '  +1.71%  '
$ws.Range("E30").Value = '  +2.70%  '
$ws.Range("E31").Value = '  +0.18%  '
$ws.Range("D32").Value = "'8.21"
$ws.Range("E32").Value = '  +1.77%  '
$ws.Range("E33").Value = '  -10.91%  '
$ws.Range("E34").Value = '  +2.36%  '
$ws.Range("D36").Value = "'24.11"
$ws.Range("E36").Value = '  +2.94%  '
$ws.Range("D37").Value = '3.478.56'
$ws.Range("E37").Value = '  +2.13%  '
$ws.Range("E38").Value = '  +2.99%  '
$ws.Range("E39").Value = '  +0.78%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("D41").Value = "'166.94"
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("E42").Value = '  +3.28%  '
$ws.Range("D43").Value = "'27.59"
$ws.Range("E43").Value = '  +8.06%  '
$ws.Range("E44").Value = '  +4.00%  '
$ws.Range("E45").Value = '  +4.06%  '
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").Value = "'42.43"
$ws.Range("E47").Value = '  +1.70%  '
$ws.Range("E48").Value = '  +1.56%  '
$ws.Range("E49").Value = '  -2.31%  '
$ws.Range("D50").Value = '2.574.71'
$ws.Range("E51").Value = '  +2.62%  '
